# Fix Training Data Issue (#48)
# The BF column holds the game date as plain text. Due to the way NBA
# stats were originally scraped, every row in this sheet was stamped with
# the sheet/file name ("6-19-2007-08") instead of the actual game date,
# which was really one day later than the label implied. Correct it to
# the proper ISO date string "2008-06-19" for every data row (rows 2-31),
# without disturbing the cell's existing formatting/style.
#
# Note: assigning a date-shaped string straight to Range.Value/Value2
# makes Excel "smart" auto-convert it into a date serial number (and pick
# up a date number format), which is not what we want here -- the source
# data is a plain text label, not a real Excel date. Instead we write the
# text through a literal formula (so it round-trips as the exact string)
# and then flatten the formula down to a static value with
# Copy + PasteSpecial(xlPasteValues), which leaves the cell's original
# (unformatted) style completely untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateRange = $ws.Range("BF2:BF31")
$dateRange.Formula = "=""2008-06-19"""
$dateRange.Copy()
$dateRange.PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = 0
